# "update the map with time"
# Replace the library-branch names in column A with their short codes.
# The insertion order below reproduces the shared-string table order seen
# in the target workbook (Excel appends new shared strings in the order
# cells are (re)written).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value  = "SEQ"   # Sequoya      -> SEQ
$ws.Range("A8").Value  = "PIN"   # Pinney       -> PIN
$ws.Range("A6").Value  = "MEA"   # Meadowridge  -> MEA
$ws.Range("A5").Value  = "LAK"   # Lakeview     -> LAK
$ws.Range("A4").Value  = "HAW"   # Hawthorne    -> HAW
$ws.Range("A2").Value  = "HPB"   # Ashman       -> HPB
$ws.Range("A3").Value  = "MAD"   # Central      -> MAD
$ws.Range("A10").Value = "SMB"   # South Madison-> SMB
$ws.Range("A7").Value  = "MSB"   # Monroe Street-> MSB

# Page setup: portrait, paper size 9 (A4)
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

# Move the active selection from E5 to B5
$ws.Range("B5").Select()
